$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.508064516129032
    "C2" = 0.513870541611625
    "D2" = 0.642725598526704
    "E2" = 0.498447204968944
    "F2" = 0.410896708286039

    "B3" = 0.6
    "C3" = 0.652575957727873
    "D3" = 0.692449355432781
    "E3" = 0.647515527950311
    "F3" = 0.52894438138479

    "B4" = 0.541935483870968
    "C4" = 0.597093791281374
    "D4" = 0.622467771639042
    "E4" = 0.545031055900621
    "F4" = 0.445327279606508
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
